$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Title value (row 5, column B) which was previously empty
$ws.Range("B5").Value = "FonctionQualifiee"

# Update the Date value (row 8, column B) to reflect the new commit date
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"
